# XAlpha.xlsx test-data fix
#
# - XAlphaDeals: fix the TestCaseID typo "X-Alpha" -> "XAlpha" in A2:A4 so the
#   naming matches QA_TestCase_Auto_XAlpha_001..003 already on XAlphaLogin.
# - XAlphaLogin: add a thin-bordered, wrap-text "SkipAtStepNum" style value
#   cell (F2:F4) matching column G's validation styling, and widen column B
#   to fit.
# - Make XAlphaDeals the active/selected sheet with a fresh selection, and
#   update XAlphaLogin's lingering selection now that it is no longer active.

$wb = $excel.ActiveWorkbook

$wsLogin = $wb.Worksheets.Item("XAlphaLogin")
$wsDeals = $wb.Worksheets.Item("XAlphaDeals")

# --- XAlphaDeals: correct the TestCaseID values (A2:A4) ---------------------
$wsDeals.Range("A2").Value = "QA_TestCase_Auto_XAlpha_004"
$wsDeals.Range("A3").Value = "QA_TestCase_Auto_XAlpha_005"
$wsDeals.Range("A4").Value = "QA_TestCase_Auto_XAlpha_006"

# --- XAlphaLogin: widen column B slightly (bestFit) --------------------------
$wsLogin.Columns.Item(2).ColumnWidth = 27.33

# --- XAlphaLogin: new thin-border / wrap-text cells F2:F4 --------------------
$fCells = $wsLogin.Range("F2:F4")
$fCells.Borders.LineStyle = 1
$fCells.Borders.Weight = 2
$fCells.WrapText = $true

# --- Selections / active sheet ------------------------------------------------
$wsLogin.Activate()
$wsLogin.Range("E10").Select()

$wsDeals.Activate()
$wsDeals.Range("A9").Select()
